$d = $word.ActiveDocument

# First occurrence: <id>p108r_1</id>  (currently split across 3 runs:
# "<id>", "p108r_1", "</id>") -> merge into a single run of text.
$d.Content.Find.Execute("<id>p108r_1</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p108r_1</id>", 2)

# Second occurrence: <id>p108r_2</id> (currently split across 4 runs:
# "<id>", "p108r_", "2", "</id>") -> merge into a single run of text.
$d.Content.Find.Execute("<id>p108r_2</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p108r_2</id>", 2)
